$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 727.72974
$ws.Range("J17").Value = 698.05884
$ws.Range("L17").Value = 2094.17652
$ws.Range("N17").Value = -2430.17652
$ws.Range("H55").Value = 361.7143
$ws.Range("I55").Value = 268.125
$ws.Range("J55").Value = 419.30768
$ws.Range("K55").Value = 268.125
$ws.Range("L55").Value = 419.30768
$ws.Range("M55").Value = -54.125
$ws.Range("N55").Value = -847.30768
$ws.Range("H111").Value = 2766.7
$ws.Range("I111").Value = 1328.1666
$ws.Range("J111").Value = 4924.5
$ws.Range("K111").Value = 3984.4998
$ws.Range("L111").Value = 14773.5
$ws.Range("M111").Value = -917.4998000000001
$ws.Range("N111").Value = -20907.5
$ws.Range("H132").Value = 100094.71
$ws.Range("I132").Value = 108057.63
$ws.Range("K132").Value = 324172.89
$ws.Range("M132").Value = -321642.89
$ws.Range("H135").Value = 2656.6667
$ws.Range("I135").Value = 984
$ws.Range("J135").Value = 3493
$ws.Range("K135").Value = 8856
$ws.Range("L135").Value = 31437
$ws.Range("M135").Value = -6321
$ws.Range("N135").Value = -36507
$ws.Range("H138").Value = 1774.2565
$ws.Range("I138").Value = 1080
$ws.Range("J138").Value = 3162.7693
$ws.Range("K138").Value = 3240
$ws.Range("L138").Value = 9488.3079
$ws.Range("M138").Value = 1900
$ws.Range("N138").Value = -19768.3079

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6292385.5
$ws.Range("I32").Value = 6292385.5
$ws.Range("K32").Value = 6292385.5
$ws.Range("M32").Value = -6292098.5
$ws.Range("H45").Value = 2933.111
$ws.Range("I45").Value = 2833
$ws.Range("J45").Value = 3133.3333
$ws.Range("K45").Value = 2833
$ws.Range("L45").Value = 3133.3333
$ws.Range("M45").Value = -2456
$ws.Range("N45").Value = -3887.3333
$ws.Range("H61").Value = 5556988
$ws.Range("I61").Value = 5556988
$ws.Range("K61").Value = 5556988
$ws.Range("M61").Value = -5556776
$ws.Range("H122").Value = 1385.7693
$ws.Range("I122").Value = 1216.8572
$ws.Range("K122").Value = 3650.5716
$ws.Range("M122").Value = -1200.5716
$ws.Range("H136").Value = 5556988
$ws.Range("I136").Value = 5556988
$ws.Range("K136").Value = 16670964
$ws.Range("M136").Value = -16668414

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3084
$ws.Range("I20").Value = 1488
$ws.Range("J20").Value = 3616
$ws.Range("K20").Value = 1488
$ws.Range("L20").Value = 3616
$ws.Range("M20").Value = -1241
$ws.Range("N20").Value = -4110
$ws.Range("H64").Value = 2001
$ws.Range("J64").Value = 2001
$ws.Range("L64").Value = 2001
$ws.Range("N64").Value = -2451
$ws.Range("H67").Value = 2001
$ws.Range("J67").Value = 2001
$ws.Range("L67").Value = 2001
$ws.Range("N67").Value = -3561
$ws.Range("H80").Value = 5812.6875
$ws.Range("I80").Value = 10688.6
$ws.Range("J80").Value = 3596.3635
$ws.Range("K80").Value = 10688.6
$ws.Range("L80").Value = 3596.3635
$ws.Range("M80").Value = -9690.6
$ws.Range("N80").Value = -5592.363499999999
$ws.Range("H83").Value = 5812.6875
$ws.Range("I83").Value = 10688.6
$ws.Range("J83").Value = 3596.3635
$ws.Range("K83").Value = 53443
$ws.Range("L83").Value = 17981.8175
$ws.Range("M83").Value = -48451
$ws.Range("N83").Value = -27965.8175
$ws.Range("H86").Value = 5166.8335
$ws.Range("J86").Value = 5006.5
$ws.Range("L86").Value = 5006.5
$ws.Range("N86").Value = -7252.5
$ws.Range("H89").Value = 5166.8335
$ws.Range("J89").Value = 5006.5
$ws.Range("L89").Value = 25032.5
$ws.Range("N89").Value = -36264.5
$ws.Range("H94").Value = 1102.6364
$ws.Range("I94").Value = 941
$ws.Range("K94").Value = 941
$ws.Range("M94").Value = -490
$ws.Range("H99").Value = 38607.47
$ws.Range("I99").Value = 48919.727
$ws.Range("J99").Value = 19701.666
$ws.Range("K99").Value = 48919.727
$ws.Range("L99").Value = 19701.666
$ws.Range("M99").Value = -47421.727
$ws.Range("N99").Value = -22697.666
$ws.Range("H134").Value = 1895624.4
$ws.Range("I134").Value = 1834993.1
$ws.Range("J134").Value = 2289728
$ws.Range("K134").Value = 5504979.300000001
$ws.Range("L134").Value = 6869184
$ws.Range("M134").Value = -5502444.300000001
$ws.Range("N134").Value = -6874254

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 210.13333
$ws.Range("I7").Value = 110.55
$ws.Range("J7").Value = 409.3
$ws.Range("K7").Value = 110.55
$ws.Range("L7").Value = 409.3
$ws.Range("M7").Value = 2.450000000000003
$ws.Range("N7").Value = -635.3
$ws.Range("H58").Value = 2062882.5
$ws.Range("I58").Value = 3086823.8
$ws.Range("K58").Value = 3086823.8
$ws.Range("M58").Value = -3086620.8
$ws.Range("H86").Value = 7863.5
$ws.Range("I86").Value = 7257.5
$ws.Range("J86").Value = 8267.5
$ws.Range("K86").Value = 7257.5
$ws.Range("L86").Value = 8267.5
$ws.Range("M86").Value = -6134.5
$ws.Range("N86").Value = -10513.5
$ws.Range("H89").Value = 7863.5
$ws.Range("I89").Value = 7257.5
$ws.Range("J89").Value = 8267.5
$ws.Range("K89").Value = 36287.5
$ws.Range("L89").Value = 41337.5
$ws.Range("M89").Value = -30671.5
$ws.Range("N89").Value = -52569.5
$ws.Range("H94").Value = 2121.818
$ws.Range("I94").Value = 1006
$ws.Range("J94").Value = 2369.7778
$ws.Range("K94").Value = 1006
$ws.Range("L94").Value = 2369.7778
$ws.Range("M94").Value = -555
$ws.Range("N94").Value = -3271.7778
$ws.Range("H105").Value = 26740.357
$ws.Range("I105").Value = 30405.5
$ws.Range("K105").Value = 30405.5
$ws.Range("M105").Value = -28658.5
$ws.Range("H132").Value = 40207772
$ws.Range("I132").Value = 47621250
$ws.Range("K132").Value = 142863750
$ws.Range("M132").Value = -142861220
$ws.Range("H134").Value = 6734.3213
$ws.Range("I134").Value = 6872.222
$ws.Range("K134").Value = 20616.666
$ws.Range("M134").Value = -18081.666
$ws.Range("H136").Value = 2062882.5
$ws.Range("I136").Value = 3086823.8
$ws.Range("K136").Value = 9260471.399999999
$ws.Range("M136").Value = -9257921.399999999

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 130.75
$ws.Range("I29").Value = 141
$ws.Range("K29").Value = 423
$ws.Range("M29").Value = -146
$ws.Range("H109").Value = 2351
$ws.Range("I109").Value = 1791.4
$ws.Range("K109").Value = 5374.200000000001
$ws.Range("M109").Value = -4334.200000000001
$ws.Range("H131").Value = 18347.924
$ws.Range("I131").Value = 1021.125
$ws.Range("J131").Value = 46070.8
$ws.Range("K131").Value = 3063.375
$ws.Range("L131").Value = 138212.4
$ws.Range("M131").Value = 1976.625
$ws.Range("N131").Value = -148292.4
$ws.Range("H137").Value = 1520.2222
$ws.Range("J137").Value = 1661.6666
$ws.Range("L137").Value = 4984.9998
$ws.Range("N137").Value = -15184.9998

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5381.7144
$ws.Range("I70").Value = 5363.8335
$ws.Range("J70").Value = 5489
$ws.Range("K70").Value = 5363.8335
$ws.Range("L70").Value = 5489
$ws.Range("M70").Value = -5093.8335
$ws.Range("N70").Value = -6029
$ws.Range("H73").Value = 5381.7144
$ws.Range("I73").Value = 5363.8335
$ws.Range("J73").Value = 5489
$ws.Range("K73").Value = 5363.8335
$ws.Range("L73").Value = 5489
$ws.Range("M73").Value = -4427.8335
$ws.Range("N73").Value = -7361
$ws.Range("H97").Value = 1969.619
$ws.Range("I97").Value = 1647
$ws.Range("J97").Value = 3905.3333
$ws.Range("K97").Value = 1647
$ws.Range("L97").Value = 3905.3333
$ws.Range("M97").Value = -1151
$ws.Range("N97").Value = -4897.3333
$ws.Range("H122").Value = 65879.766
$ws.Range("I122").Value = 94077.82000000001
$ws.Range("K122").Value = 282233.46
$ws.Range("M122").Value = -279783.46
$ws.Range("H135").Value = 94952
$ws.Range("J135").Value = 94952
$ws.Range("L135").Value = 94952
$ws.Range("N135").Value = -105092

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 9817.125
$ws.Range("I46").Value = 23880.334
$ws.Range("J46").Value = 1379.2
$ws.Range("K46").Value = 23880.334
$ws.Range("L46").Value = 1379.2
$ws.Range("M46").Value = -23692.334
$ws.Range("N46").Value = -1755.2
$ws.Range("H68").Value = 3938.5
$ws.Range("I68").Value = 6500
$ws.Range("J68").Value = 3084.6667
$ws.Range("K68").Value = 6500
$ws.Range("L68").Value = 3084.6667
$ws.Range("M68").Value = -5751
$ws.Range("N68").Value = -4582.6667
$ws.Range("H71").Value = 3938.5
$ws.Range("I71").Value = 6500
$ws.Range("J71").Value = 3084.6667
$ws.Range("K71").Value = 32500
$ws.Range("L71").Value = 15423.3335
$ws.Range("M71").Value = -28756
$ws.Range("N71").Value = -22911.3335
$ws.Range("H122").Value = 4999.8335
$ws.Range("I122").Value = 4749.8125
$ws.Range("K122").Value = 14249.4375
$ws.Range("M122").Value = -11799.4375
$ws.Range("H136").Value = 71412.72
$ws.Range("I136").Value = 2390
$ws.Range("J136").Value = 157691.12
$ws.Range("K136").Value = 7170
$ws.Range("L136").Value = 473073.36
$ws.Range("M136").Value = -4620
$ws.Range("N136").Value = -478173.36
$ws.Range("H140").Value = 120429
$ws.Range("J140").Value = 120429
$ws.Range("L140").Value = 120429
$ws.Range("N140").Value = -130789

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1626.5
$ws.Range("I100").Value = 1305.3334
$ws.Range("J100").Value = 1947.6666
$ws.Range("K100").Value = 2610.6668
$ws.Range("L100").Value = 3895.3332
$ws.Range("M100").Value = -2069.6668
$ws.Range("N100").Value = -4977.3332
$ws.Range("H122").Value = 2610.6
$ws.Range("I122").Value = 2134.5
$ws.Range("K122").Value = 6403.5
$ws.Range("M122").Value = -3953.5
$ws.Range("H126").Value = 6720.6924
$ws.Range("I126").Value = 6596.6665
$ws.Range("K126").Value = 19789.9995
$ws.Range("M126").Value = -17319.9995
$ws.Range("H136").Value = 15153.5
$ws.Range("I136").Value = 17750.5
$ws.Range("K136").Value = 53251.5
$ws.Range("M136").Value = -50701.5
